$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.162.16"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "2.053.50"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'248.47"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = "  -2.11%  "

$ws.Range("D7").Value = "'57.50"
$ws.Range("E7").Value = "  -2.81%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.381"
$ws.Range("E9").Value = "  -2.63%  "

$ws.Range("D10").Value = "'0.0779"
$ws.Range("E10").Value = "  -2.90%  "

$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("D12").Value = "'16.03"
$ws.Range("E12").Value = "  -1.82%  "

$ws.Range("D13").Value = "'0.888"
$ws.Range("E13").Value = "  +7.81%  "

$ws.Range("D14").Value = "2.353.26"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("E15").Value = "  +3.28%  "

$ws.Range("D16").Value = "2.053.31"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "'18.19"
$ws.Range("E17").Value = "  +15.02%  "

$ws.Range("D18").Value = "37.188.05"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").Value = "'74.76"
$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").Value = "0.0₃0893"
$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("D21").Value = "'5.39"
$ws.Range("E21").Value = "  -1.69%  "

$ws.Range("D22").Value = "'237.32"
$ws.Range("E22").Value = "  -0.97%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("E24").Value = "  +2.92%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.52"
$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  -3.76%  "

$ws.Range("D27").Value = "'169.58"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").Value = "'20.08"
$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.14"
$ws.Range("E30").Value = "  -1.09%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.86"
$ws.Range("E31").Value = "  +1.21%  "

$ws.Range("D32").Value = "'0.0620"
$ws.Range("E32").Value = "  -2.21%  "

$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("D34").Value = "'0.0892"
$ws.Range("E34").Value = "  -1.98%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  -2.11%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("E38").Value = "  -1.97%  "

$ws.Range("D39").Value = "'5.29"
$ws.Range("E39").Value = "  +14.59%  "

$ws.Range("D40").Value = "'3.13"
$ws.Range("E40").Value = "  +10.52%  "

$ws.Range("D41").Value = "'0.0989"
$ws.Range("E41").Value = "  -14.66%  "

$ws.Range("E42").Value = "  -1.85%  "

$ws.Range("D43").Value = "'17.28"
$ws.Range("E43").Value = "  -3.64%  "

$ws.Range("D44").Value = "'1.14"
$ws.Range("E44").Value = "  -2.56%  "

$ws.Range("D45").Value = "'96.11"
$ws.Range("E45").Value = "  -3.11%  "

$ws.Range("D46").Value = "'2.44"
$ws.Range("E46").Value = "  -2.07%  "

$ws.Range("D47").Value = "1.269.53"
$ws.Range("E47").Value = "  -2.88%  "

$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -3.03%  "

$ws.Range("D49").Value = "'6.82"
$ws.Range("E49").Value = "  -1.70%  "

$ws.Range("D50").Value = "2.241.03"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("D51").Value = "'43.95"
$ws.Range("E51").Value = "  -0.91%  "
